# Longitud_Vias_Ferreas.xlsx - "Add files via upload" update
# Refreshes the 2023 railway-length table to 2024 figures (ARTF source),
# updating the title/footer captions and the changed state totals.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("C_20.1")

# --- Title (B2) and footer update line (B37) -----------------------------
$ws.Range("B2").Value  = "Longitud de vías férreas 2024"
$ws.Range("B37").Value = "Actualización: Julio 2025."

# --- Updated state figures (Troncales/ramales + Secundarias columns) -----
# Campeche
$ws.Range("C8").Value  = 875.94977350477279
$ws.Range("D8").Value  = 52.299000000000007

# Ciudad de México
$ws.Range("C13").Value = 110.30766992188222
$ws.Range("D13").Value = 101.482

# México
$ws.Range("C19").Value = 924.8953312834218
$ws.Range("D19").Value = 262.18200000000002
$ws.Range("E19").Value = 89.72596999999999

# Quintana Roo
$ws.Range("C27").Value = 966.55500000000006
$ws.Range("D27").Value = 47.957999999999998

# Tabasco
$ws.Range("D31").Value = 13.654

# Yucatán
$ws.Range("C35").Value = 1058.538861934848
$ws.Range("D35").Value = 34.557000000000002

# --- Minor layout refresh: title row shrinks along with the new default --
# row height used when the workbook was re-saved from the newer Excel build.
$ws.Rows.Item(2).RowHeight = 18

Write-Host "Longitud_Vias_Ferreas: updated to 2024 figures"
